# Weekly price update: insert a new data row (week of 2022-05-26) at the
# top of the price history table for "Coco" and push the existing rows
# (previously rows 17-65) down by one, turning the last row (65) into 66.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 17, shifting rows 17..65 down to 18..66.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with this week's record. All the
# descriptive columns (market/region/product taxonomy, unit, origin,
# kg-per-unit) are constant across every row in this sheet, so reuse the
# same values; only the date, volume and price columns are new.
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44707
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100108
$ws.Range("H17").Value = "Tropicales y subtropicales"
$ws.Range("I17").Value = 100108007
$ws.Range("J17").Value = "Coco"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 25
$ws.Range("N17").Value = 30000
$ws.Range("O17").Value = 30000
$ws.Range("P17").Value = 30000
$ws.Range("Q17").Value = "$/malla 20 unidades"
$ws.Range("R17").Value = "Perú"
$ws.Range("S17").Value = 1500
$ws.Range("T17").Value = 20
